$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-21 01:40:00"

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
